$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ten_lists")

# Outer list markings (column C/I) - these relabel which sub-list is which,
# plus new "inner" markers in columns E/K ("start SD"/"start DD") next to
# each outer label, and a new column M flag ("didn't do this trial") on the
# rows that start the DD (no-walk) trials.

$ws.Range("C3").Value = "no walk/same"
$ws.Range("E3").Value = "start SD"
$ws.Range("I3").Value = "no walk/diff"
$ws.Range("K3").Value = "start SD"

$ws.Range("C10").Value = "walk/diff"
$ws.Range("E10").Value = "start DD"
$ws.Range("I10").Value = "no walk/same"
$ws.Range("K10").Value = "start DD"

$ws.Range("C17").Value = "walk/same"
$ws.Range("E17").Value = "start SD"
$ws.Range("I17").Value = "walk/diff"
$ws.Range("K17").Value = "start SD"

$ws.Range("C24").Value = "no walk/diff"
$ws.Range("E24").Value = "start DD"
$ws.Range("I24").Value = "walk/same"
$ws.Range("K24").Value = "start SD"

$ws.Range("C31").Value = "walk/same"
$ws.Range("E31").Value = "start DD"
$ws.Range("I31").Value = "no walk/diff"
$ws.Range("K31").Value = "start DD"

$ws.Range("M25").Value = "didn't do this trial"
$ws.Range("M32").Value = "didn't do this trial"

$ws.Range("M33").Select()
